$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking data refresh (prices / 1h volume %, plus two rows whose
# rank order swapped: RenzoRestakedETH<->Aptos at 37/38 and
# Monero<->OKB at 47/48).
#
# Force the whole data block to Text format first so that numeric-
# looking values (e.g. "165.60", "0.520") are stored as literal text
# -- matching the original inlineStr cells -- rather than being
# auto-coerced into Number cells (which would silently drop
# significant trailing zeros). The style is restored to Normal at the
# end so no visible formatting change is left behind.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# --- Simple price / volume(1h) updates ---
$ws.Range("D2").Value = "69.341.19"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "3.782.54"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").Value = "624.01"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").Value = "165.60"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("D7").Value = "3.781.89"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").Value = "6.69"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "35.66"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "4.415.99"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "3.781.00"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "69.298.27"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "17.67"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "7.12"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D21").Value = "468.47"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").Value = "9.64"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").Value = "0.704"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +4.96%  "
$ws.Range("D25").Value = "83.32"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "3.930.78"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +3.40%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").Value = "7.28"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D34").Value = "28.84"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +17.17%  "
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  +8.38%  "
$ws.Range("D41").Value = "5.82"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").Value = "0.970"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "43.43"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D49").Value = "1.92"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").Value = "  +0.37%  "

# --- Row 37/38 swap (RenzoRestakedETH <-> Aptos) with updated values ---
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "9.01"
$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.731.89"
$ws.Range("E38").Value = "  +0.39%  "

# --- Row 47/48 swap (Monero <-> OKB) with updated values ---
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "46.78"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "152.21"
$ws.Range("E48").Value = "  +0.00%  "

# Restore default (General) style on the data block now that all
# values are safely stored as text.
$dataRange.Style = "Normal"
